$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 ("Lode Runner sur une barre de franchissement") ---
# Responsable reassigned from Ankit to Ribensky
$ws.Range("E14").Value = "Ribensky"

# --- Row 15 ("Lode Runner qui tombe") : task completed ---
# Pick up the "completed" date-cell fill (used elsewhere, e.g. C9/D9) for C15:D15
$ws.Range("C9").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)
# Responsable swaps back from Ribensky to Ankit
$ws.Range("E15").Value = "Ankit"
# Progression becomes "Fini", with the "completed" progression fill (used e.g. at F6)
$ws.Range("F6").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "Fini"

# --- Row 27 ("Pointage") : task completed ---
$ws.Range("C9").Copy()
$ws.Range("C27:D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "Ankit"
$ws.Range("F6").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = "Fini"

# --- Row 28 ("Prochain niveau") : task completed ---
$ws.Range("C9").Copy()
$ws.Range("C28:D28").PasteSpecial(-4122)
$ws.Range("F6").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = "Fini"

# --- Update the active selection to reflect where the user left off ---
$ws.Range("H14").Select() | Out-Null

$excel.CutCopyMode = 0
